$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.418.04'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = '3.476.28'
$ws.Range('E3').Value = '  -4.33%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'578.47"
$ws.Range('E5').Value = '  -4.42%  '
$ws.Range('D6').Value = "'192.13"
$ws.Range('E6').Value = '  -3.56%  '
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('D8').Value = '3.461.80'
$ws.Range('E8').Value = '  -4.39%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('E10').Value = '  -7.26%  '
$ws.Range('E11').Value = '  -4.76%  '
$ws.Range('D12').Value = "'51.35"
$ws.Range('E12').Value = '  -4.83%  '
$ws.Range('E13').Value = '  -6.64%  '
$ws.Range('D14').Value = "'9.14"
$ws.Range('E14').Value = '  -4.48%  '
$ws.Range('D15').Value = '4.033.15'
$ws.Range('E15').Value = '  -4.33%  '
$ws.Range('D16').Value = "'653.44"
$ws.Range('E16').Value = '  -3.72%  '
$ws.Range('D17').Value = '69.302.64'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '3.476.89'
$ws.Range('E18').Value = '  -5.42%  '
$ws.Range('E19').Value = '  -5.26%  '
$ws.Range('E20').Value = '  -1.73%  '
$ws.Range('D21').Value = "'18.19"
$ws.Range('E21').Value = '  -4.73%  '
$ws.Range('D22').Value = "'0.944"
$ws.Range('E22').Value = '  -5.58%  '
$ws.Range('D23').Value = "'18.06"
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').Value = "'98.98"
$ws.Range('E25').Value = '  -6.61%  '
$ws.Range('E26').Value = '  -7.60%  '
$ws.Range('E27').Value = '  -4.34%  '
$ws.Range('D28').Value = "'9.99"
$ws.Range('E28').Value = '  -4.34%  '
$ws.Range('E29').Value = '  -4.92%  '
$ws.Range('D30').Value = "'32.50"
$ws.Range('E30').Value = '  -4.51%  '
$ws.Range('E31').Value = '  -8.90%  '
$ws.Range('D32').Value = "'6.73"
$ws.Range('E32').Value = '  -6.71%  '
$ws.Range('E33').Value = '  -5.15%  '
$ws.Range('E34').Value = '  -5.57%  '
$ws.Range('D35').Value = "'60.82"
$ws.Range('E35').Value = '  -4.17%  '
$ws.Range('D36').Value = '3.720.67'
$ws.Range('E36').Value = '  -6.73%  '
$ws.Range('D37').Value = "'525.29"
$ws.Range('E37').Value = '  +3.52%  '
$ws.Range('D38').Value = "'1.00"
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').Value = '0.0₃0791'
$ws.Range('E39').Value = '  -8.69%  '
$ws.Range('D40').Value = "'2.92"
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('D41').Value = "'3.49"
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('E42').Value = '  -4.10%  '
$ws.Range('D43').Value = "'0.133"
$ws.Range('E43').Value = '  -2.39%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = "'34.26"
$ws.Range('E44').Value = '  -7.06%  '
$ws.Range('B45').Value = 'CoreDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D45').Value = "'3.50"
$ws.Range('E45').Value = '  +71.29%  '
$ws.Range('D46').Value = "'0.0442"
$ws.Range('E46').Value = '  -4.28%  '
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('D48').Value = "'2.83"
$ws.Range('E48').Value = '  -9.20%  '
$ws.Range('E49').Value = '  -4.76%  '
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').Value = "'8.15"
$ws.Range('E51').Value = '  -6.14%  '
